$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("2019 South")

$oldA = $ws4.Range("A14").Value()
$oldB = $ws4.Range("B14").Value()
$oldD = $ws4.Range("D14").Value()
$oldF = $ws4.Range("F14").Value()
$oldH = $ws4.Range("H14").Value()
$oldJ = $ws4.Range("J14").Value()
$oldL = $ws4.Range("L14").Value()
$oldN = $ws4.Range("N14").Value()
$oldP = $ws4.Range("P14").Value()
$oldR = $ws4.Range("R14").Value()
$oldT = $ws4.Range("T14").Value()
$oldV = $ws4.Range("V14").Value()
$oldX = $ws4.Range("X14").Value()
$oldZ = $ws4.Range("Z14").Value()
$oldAB = $ws4.Range("AB14").Value()
$oldAD = $ws4.Range("AD14").Value()
$oldAF = $ws4.Range("AF14").Value()
$oldAH = $ws4.Range("AH14").Value()
$oldAJ = $ws4.Range("AJ14").Value()

$ws4.Rows.Item(14).Delete()
$ws4.Rows.Item(18).Insert()

$ws4.Range("A18").Value = $oldA
$ws4.Range("B18").Value = "Riada, Liadh Ninot-elected"
$ws4.Range("D18").Value = $oldD
$ws4.Range("E18").Formula = "=IF(ISBLANK(F18),-D18,F18-D18)"
$ws4.Range("F18").Value = $oldF
$ws4.Range("G18").Formula = "=IF(ISBLANK(H18),-F18,H18-F18)"
$ws4.Range("H18").Value = $oldH
$ws4.Range("I18").Formula = "=IF(ISBLANK(J18),-H18,J18-H18)"
$ws4.Range("J18").Value = $oldJ
$ws4.Range("K18").Formula = "=IF(ISBLANK(L18),-J18,L18-J18)"
$ws4.Range("L18").Value = $oldL
$ws4.Range("M18").Formula = "=IF(ISBLANK(N18),-L18,N18-L18)"
$ws4.Range("N18").Value = $oldN
$ws4.Range("O18").Formula = "=IF(ISBLANK(P18),-N18,P18-N18)"
$ws4.Range("P18").Value = $oldP
$ws4.Range("Q18").Formula = "=IF(ISBLANK(R18),-P18,R18-P18)"
$ws4.Range("R18").Value = $oldR
$ws4.Range("S18").Formula = "=IF(ISBLANK(T18),-R18,T18-R18)"
$ws4.Range("T18").Value = $oldT
$ws4.Range("U18").Formula = "=IF(ISBLANK(V18),-T18,V18-T18)"
$ws4.Range("V18").Value = $oldV
$ws4.Range("W18").Formula = "=IF(ISBLANK(X18),-V18,X18-V18)"
$ws4.Range("X18").Value = $oldX
$ws4.Range("Y18").Formula = "=IF(ISBLANK(Z18),-X18,Z18-X18)"
$ws4.Range("Z18").Value = $oldZ
$ws4.Range("AA18").Formula = "=IF(ISBLANK(AB18),-Z18,AB18-Z18)"
$ws4.Range("AB18").Value = $oldAB
$ws4.Range("AC18").Formula = "=IF(ISBLANK(AD18),-AB18,AD18-AB18)"
$ws4.Range("AD18").Value = $oldAD
$ws4.Range("AE18").Formula = "=IF(ISBLANK(AF18),-AD18,AF18-AD18)"
$ws4.Range("AF18").Value = $oldAF
$ws4.Range("AG18").Formula = "=IF(ISBLANK(AH18),-AF18,AH18-AF18)"
$ws4.Range("AH18").Value = $oldAH
$ws4.Range("AI18").Formula = "=IF(ISBLANK(AJ18),-AH18,AJ18-AH18)"
$ws4.Range("AJ18").Value = $oldAJ
Write-Output "done"
